# Regenerate merged AHB files
# - rename header row columns: "_old" -> "_FV2404", "_new" -> "_FV2410"
# - wrap the data range in a table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels (row 1) ---------------------------------
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_old$", "_FV2404")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Text -replace "_new$", "_FV2410")
}

# --- 2. Turn the used range into a table (ListObject) ---------------------
$rng = $ws.Range("A1:U80")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
